$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 35
$ws.Range("H3").Value = 36
$ws.Range("E4").Value = 50
$ws.Range("E5").Value = 156
$ws.Range("F5").Value = 109
$ws.Range("H5").Value = 120
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 36
$ws.Range("H6").Value = 46
$ws.Range("E10").Value = 675
$ws.Range("F10").Value = 376
$ws.Range("G10").Value = 95
$ws.Range("H10").Value = 471
$ws.Range("E11").Value = 443
$ws.Range("G11").Value = 65
$ws.Range("H11").Value = 313
$ws.Range("E12").Value = 680
$ws.Range("F12").Value = 412
$ws.Range("H12").Value = 498
$ws.Range("E13").Value = 162
$ws.Range("F13").Value = 92
$ws.Range("H13").Value = 126
$ws.Range("E14").Value = 140
$ws.Range("F14").Value = 82
$ws.Range("H14").Value = 116
$ws.Range("E15").Value = 197
$ws.Range("F15").Value = 94
$ws.Range("H15").Value = 145
$ws.Range("E16").Value = 231
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 32
$ws.Range("H18").Value = 49
$ws.Range("E20").Value = 99
$ws.Range("F20").Value = 44
$ws.Range("H20").Value = 81
$ws.Range("E22").Value = 194
$ws.Range("F22").Value = 108
$ws.Range("H22").Value = 150
$ws.Range("E23").Value = 223
$ws.Range("F23").Value = 115
$ws.Range("H23").Value = 167
$ws.Range("E24").Value = 258
$ws.Range("F24").Value = 152
$ws.Range("H24").Value = 182
$ws.Range("E25").Value = 322
$ws.Range("F25").Value = 179
$ws.Range("H25").Value = 239
$ws.Range("E26").Value = 186
$ws.Range("F26").Value = 117
$ws.Range("H26").Value = 142
$ws.Range("E27").Value = 377
$ws.Range("F27").Value = 206
$ws.Range("H27").Value = 288
$ws.Range("E28").Value = 222
$ws.Range("E29").Value = 192
$ws.Range("E30").Value = 249
$ws.Range("E32").Value = 210
$ws.Range("F32").Value = 134
$ws.Range("H32").Value = 172
$ws.Range("E33").Value = 323
$ws.Range("E34").Value = 245
$ws.Range("F34").Value = 174
$ws.Range("H34").Value = 212
$ws.Range("E35").Value = 178
$ws.Range("F35").Value = 125
$ws.Range("H35").Value = 152
$ws.Range("E36").Value = 90
$ws.Range("F36").Value = 59
$ws.Range("H36").Value = 69
$ws.Range("E37").Value = 190
$ws.Range("E38").Value = 105
$ws.Range("E39").Value = 194
$ws.Range("E40").Value = 302
$ws.Range("F40").Value = 156
$ws.Range("H40").Value = 236
$ws.Range("E41").Value = 431
$ws.Range("F41").Value = 217
$ws.Range("H41").Value = 309
$ws.Range("E42").Value = 448
$ws.Range("F42").Value = 257
$ws.Range("H42").Value = 318
$ws.Range("E44").Value = 355
$ws.Range("E46").Value = 384
$ws.Range("F46").Value = 224
$ws.Range("H46").Value = 288
$ws.Range("E47").Value = 531
$ws.Range("F47").Value = 299
$ws.Range("H47").Value = 391
$ws.Range("E48").Value = 261
$ws.Range("F48").Value = 127
$ws.Range("H48").Value = 171
$ws.Range("E49").Value = 335
$ws.Range("F49").Value = 168
$ws.Range("H49").Value = 255
$ws.Range("E50").Value = 273
$ws.Range("F50").Value = 149
$ws.Range("H50").Value = 222
$ws.Range("E51").Value = 265
$ws.Range("E52").Value = 31
$ws.Range("F52").Value = 14
$ws.Range("H52").Value = 22
